# Adds the two new match rows (125 and 126) to the Ekstraklasa sheet,
# mirroring the styling of the last existing data row (124):
#   - column A uses the bold/bordered/centered style
#   - column E uses the date-time number format
# then fills in all the data from the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the formatting (style indices) of the last data row down into
# the two new rows before writing values into them.
$ws.Range("A124:V124").Copy()
$ws.Range("A125:V126").PasteSpecial(-4122)
$ws.Range("A1").Select()

$newRows = @(
    @{
        Row = 125
        A = 124
        B = "poland"
        C = "ekstraklasa"
        D = "2023-2024"
        E = 45240.75
        F = "Warta Poznan"
        G = 0
        H = "Puszcza"
        I = 2
        J = 1.78
        K = "04/11/2023 17:42"
        L = 2.11
        M = "10/11/2023 17:56"
        N = 3.45
        O = "04/11/2023 17:42"
        P = 3.1
        Q = "10/11/2023 17:56"
        R = 5.12
        S = "04/11/2023 17:42"
        T = 4.2
        U = "10/11/2023 17:56"
        V = "https://www.betexplorer.com/football/poland/ekstraklasa/warta-poznan-puszcza/8E3HhJpi/"
    },
    @{
        Row = 126
        A = 125
        B = "poland"
        C = "ekstraklasa"
        D = "2023-2024"
        E = 45240.85416666666
        F = "Cracovia"
        G = 0
        H = "Slask Wroclaw"
        I = 1
        J = 2.18
        K = "06/11/2023 19:12"
        L = 2.18
        M = "10/11/2023 20:29"
        N = 3.32
        O = "06/11/2023 19:12"
        P = 3.26
        Q = "10/11/2023 20:28"
        R = 3.31
        S = "06/11/2023 19:12"
        T = 3.71
        U = "10/11/2023 20:29"
        V = "https://www.betexplorer.com/football/poland/ekstraklasa/cracovia-slask-wroclaw/Y1xnwvo4/"
    }
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in $columns) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}

$wb.Save()
